$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text, so values like "1.000" or "0.9992"
# are not reinterpreted as numbers and keep their exact original formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.365.76"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.04"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.85"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3798"
$ws.Range("E7").Value = "  +3.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3283"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.32"
$ws.Range("E9").Value = "  -8.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.153"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07408"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9993"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.57"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.865"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.825"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.570.17"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001091"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06714"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.32"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.416"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.30"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.365.52"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.584"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.50"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.46"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.948"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.54"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.745.95"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.087"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.979"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.612"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08301"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02403"
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2187"
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06289"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.15"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6110"
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.91"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5956"
$ws.Range("E46").Value = "  -4.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.752"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.010"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.34"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.179"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07109"
$ws.Range("E51").Value = "  -2.40%  "

# Rows 38 and 39 swap places (TrustWalletToken <-> InternetComputer(DFINITY))
# including their price/volume values.
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.361"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.284"
$ws.Range("E39").Value = "  -1.35%  "
